$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$Barr = New-Object "object[,]" 24,1
$Barr[0,0] = 0.1424735206852858
$Barr[1,0] = 0.132922669392002
$Barr[2,0] = 0.1271330594934312
$Barr[3,0] = 0.124792503172543
$Barr[4,0] = 0.124404988218032
$Barr[5,0] = 0.1271014179985173
$Barr[6,0] = 0.1391648776732239
$Barr[7,0] = 0.1634162523099576
$Barr[8,0] = 0.1816030026128033
$Barr[9,0] = 0.1899586220687723
$Barr[10,0] = 0.1931346474524531
$Barr[11,0] = 0.1924501012305626
$Barr[12,0] = 0.1902196758566532
$Barr[13,0] = 0.188855032920145
$Barr[14,0] = 0.1810586080036671
$Barr[15,0] = 0.1762969098381575
$Barr[16,0] = 0.1735658543289134
$Barr[17,0] = 0.1726424940451636
$Barr[18,0] = 0.1768029983294781
$Barr[19,0] = 0.1908744811963601
$Barr[20,0] = 0.2001406197714175
$Barr[21,0] = 0.1951887044031508
$Barr[22,0] = 0.1765741754572758
$Barr[23,0] = 0.1567912892315917
$ws.Range("B2:B25").Value = $Barr

$Darr = New-Object "object[,]" 24,1
$Darr[0,0] = 0.1075218374974156
$Darr[1,0] = 0.09944267508923588
$Darr[2,0] = 0.09503150773295488
$Darr[3,0] = 0.09336411875966633
$Darr[4,0] = 0.09309488664497678
$Darr[5,0] = 0.0950085038802797
$Darr[6,0] = 0.1046171998833927
$Darr[7,0] = 0.1281959620443445
$Darr[8,0] = 0.148956342421144
$Darr[9,0] = 0.1592762868475575
$Darr[10,0] = 0.1633216129475272
$Darr[11,0] = 0.1624440793620749
$Darr[12,0] = 0.1596062767239061
$Darr[13,0] = 0.1578862842723652
$Darr[14,0] = 0.1483004962890391
$Darr[15,0] = 0.1426521796687723
$Darr[16,0] = 0.1394851394193495
$Darr[17,0] = 0.1384265334330621
$Darr[18,0] = 0.1432449142310759
$Darr[19,0] = 0.1604359847505634
$Darr[20,0] = 0.1724779874925559
$Darr[21,0] = 0.165973163959876
$Darr[22,0] = 0.1429766898226887
$Darr[23,0] = 0.121255340534816
$ws.Range("D2:D25").Value = $Darr

$Earr = New-Object "object[,]" 24,1
$Earr[0,0] = 0.9508387259566859
$Earr[1,0] = 0.8282990960588279
$Earr[2,0] = 0.7530571070131202
$Earr[3,0] = 0.7223868391943995
$Earr[4,0] = 0.7172933280688198
$Earr[5,0] = 0.7526435217259007
$Earr[6,0] = 0.9085824366776478
$Earr[7,0] = 1.214773457444011
$Earr[8,0] = 1.440623796129159
$Earr[9,0] = 1.543724221261954
$Earr[10,0] = 1.582831794396697
$Earr[11,0] = 1.574406122682319
$Earr[12,0] = 1.546940224353023
$Earr[13,0] = 1.530125556306587
$Earr[14,0] = 1.433894345872773
$Earr[15,0] = 1.374961327013665
$Earr[16,0] = 1.341097589705157
$Earr[17,0] = 1.329637169889111
$Earr[18,0] = 1.381231337671551
$Earr[19,0] = 1.555005724824639
$Earr[20,0] = 1.668967533220808
$Earr[21,0] = 1.608103340039975
$Earr[22,0] = 1.378396611944567
$Earr[23,0] = 1.131834713763794
$ws.Range("E2:E25").Value = $Earr

$Farr = New-Object "object[,]" 24,1
$Farr[0,0] = 2.576855933070419
$Farr[1,0] = 2.395756840799038
$Farr[2,0] = 2.29177132131386
$Farr[3,0] = 2.251114324543266
$Farr[4,0] = 2.244464312228587
$Farr[5,0] = 2.291216173841519
$Farr[6,0] = 2.512858060101564
$Farr[7,0] = 3.009204117294473
$Farr[8,0] = 3.418071009080109
$Farr[9,0] = 3.615218325157912
$Farr[10,0] = 3.691613079526007
$Farr[11,0] = 3.675080490328071
$Farr[12,0] = 3.621467700192511
$Farr[13,0] = 3.588859007953374
$Farr[14,0] = 3.405422779271476
$Farr[15,0] = 3.295840743611905
$Farr[16,0] = 3.233853093034753
$Farr[17,0] = 3.213039940541933
$Farr[18,0] = 3.307397238718551
$Farr[19,0] = 3.637166746173193
$Farr[20,0] = 3.862897757224914
$Farr[21,0] = 3.741439732502556
$Farr[22,0] = 3.302169403665857
$Farr[23,0] = 2.867665339579503
$ws.Range("F2:F25").Value = $Farr

$Garr = New-Object "object[,]" 24,1
$Garr[0,0] = 0.002375260358286178
$Garr[1,0] = 0.002393342844594032
$Garr[2,0] = 0.002404901984110808
$Garr[3,0] = 0.002409728540122307
$Garr[4,0] = 0.002410537037506991
$Garr[5,0] = 0.002404966604760998
$Garr[6,0] = 0.002381401294057062
$Garr[7,0] = 0.00233874770379932
$Garr[8,0] = 0.002309486177406484
$Garr[9,0] = 0.002296603623999105
$Garr[10,0] = 0.002291785128198869
$Garr[11,0] = 0.002292820245307073
$Garr[12,0] = 0.002296206014776558
$Garr[13,0] = 0.002298287635342544
$Garr[14,0] = 0.002310336557128673
$Garr[15,0] = 0.00231783672526873
$Garr[16,0] = 0.002322191084635539
$Garr[17,0] = 0.002323672397915146
$Garr[18,0] = 0.002317034144777328
$Garr[19,0] = 0.002295209922799281
$Garr[20,0] = 0.002281294462021718
$Garr[21,0] = 0.002288690188957878
$Garr[22,0] = 0.002317396859210234
$Garr[23,0] = 0.002349915068826353
$ws.Range("G2:G25").Value = $Garr

$Marr = New-Object "object[,]" 24,1
$Marr[0,0] = 5.476121970505687
$Marr[1,0] = 4.798433118031454
$Marr[2,0] = 4.383540310581708
$Marr[3,0] = 4.214704252232792
$Marr[4,0] = 4.186681503038528
$Marr[5,0] = 4.381262459459151
$Marr[6,0] = 5.242160787821092
$Marr[7,0] = 6.94338190416272
$Marr[8,0] = 8.206386836440799
$Marr[9,0] = 8.785068669136649
$Marr[10,0] = 9.00490653907309
$Marr[11,0] = 8.957527397084107
$Marr[12,0] = 8.803140105795137
$Marr[13,0] = 8.708668448662252
$Marr[14,0] = 8.168660944957139
$Marr[15,0] = 7.838517725208021
$Marr[16,0] = 7.649010211833172
$Marr[17,0] = 7.5849086710503
$Marr[18,0] = 7.873621722137671
$Marr[19,0] = 8.84846738730846
$Marr[20,0] = 9.489743557682459
$Marr[21,0] = 9.147063596260296
$Marr[22,0] = 7.857750290605054
$Marr[23,0] = 6.481252232435963
$ws.Range("M2:M25").Value = $Marr
